# Add 2023 (Batumi / "genders") column S to the trade worksheet,
# mirroring the formatting of the existing 2022 column (R), and update
# the selection/merge/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for each row in column R -> column S, then set the
# new values (Copy() clones number formats/styles exactly, matching
# how the original sheet built out each year column).
$ws.Range("R1").Copy($ws.Range("S1"))

$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2023

$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 1383.1

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 228.6

$ws.Range("R6").Copy($ws.Range("S6"))
$ws.Range("S6").Value = 7469

$ws.Range("R7").Copy($ws.Range("S7"))
$ws.Range("S7").Value = 5226

$ws.Range("R8").Copy($ws.Range("S8"))
$ws.Range("S8").Value = 724.1

$ws.Range("R9").Copy($ws.Range("S9"))
$ws.Range("S9").Value = 56.3

$ws.Range("R10").Copy($ws.Range("S10"))
$ws.Range("S10").Value = 46.2

$ws.Range("R11").Copy($ws.Range("S11"))
$ws.Range("S11").Value = 172.2

$ws.Range("R12").Copy($ws.Range("S12"))
$ws.Range("S12").Value = 13

$ws.Range("R13").Copy($ws.Range("S13"))
$ws.Range("S13").Value = 1190.7

$ws.Range("R14").Copy($ws.Range("S14"))
$ws.Range("S14").Value = 1145.9000000000001

# Extend the title merge from A1:R1 to A1:S1.
$ws.Range("A1:S1").Merge()

# Match the saved selection/active cell from the diff.
$ws.Range("S3:S14").Select()
